# Applies the "Append: 2025-09-30 18:22 JST" scrape update to the
# "ランサーズ" sheet:
#  - 3 new job rows are inserted (after the 2 rows that stayed put),
#    pushing the previously-existing rows down by 3.
#  - every data row's "取得日時" (column A) timestamp is refreshed to the
#    new scrape time.
#  - column D is widened slightly (28 -> 30 chars).
#  - hyperlinks on column F are (re)established for every data row so the
#    displayed URL text and the underlying link target stay in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-30 18:22:40"

# ---------------------------------------------------------------------
# 1) Widen column D (28 -> 30 characters). This engine's ColumnWidth
#    setter adds a fixed 5/6-character padding when it stores the width,
#    so back that out to land exactly on the target stored width.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 30 - 5/6

# ---------------------------------------------------------------------
# 2) Make room for the 3 new postings, directly below the two rows
#    (2 & 3) that are unaffected content-wise, i.e. before the old row 4.
# ---------------------------------------------------------------------
$ws.Rows("4:6").Insert()

# ---------------------------------------------------------------------
# 3) Full target content for every data row (2..12) after the insert.
# ---------------------------------------------------------------------
$rows = @(
  @{R=2;  B="【限定タスク】SIM AI の Google 認証ログイン機能の「最終調整」のみ代行(環境構築済み)"; C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定";       E="期限情報なし"; F="https://www.lancers.jp/work/detail/5403583"; G=298; H="🔥AI,Ai"}
  @{R=3;  B="【急募】メモリデータ管理ツール開発のプロフェッショナル募集";                                E="期限情報なし"; C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定";       F="https://www.lancers.jp/work/detail/5393508"; G=158; H="◆ツール,開発 ◇管理"}
  @{R=4;  B="システム開発において活躍できる案件紹介";                                         C="システム開発"; D="500,000 円 ~ 1,000,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5397117"; G=125; H="◆開発,システム開発"}
  @{R=5;  B="検索エンジン予測変換(サジェスト表示)に関するツール開発・調査依頼";                           C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定";       E="期限情報なし"; F="https://www.lancers.jp/work/detail/5403988"; G=123; H="◆ツール,開発"}
  @{R=6;  B="【3万円/Webツール開発】サジェスト対策";                                            C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定";       E="期限情報なし"; F="https://www.lancers.jp/work/detail/5403789"; G=123; H="◆ツール,開発"}
  @{R=7;  B="【RPA/Power Automate】税務システム自動化プロジェクトの依頼";                               C="システム開発"; D="50,000 円 ~ 100,000 円 / 固定";      E="期限情報なし"; F="https://www.lancers.jp/work/detail/5403634"; G=103; H="◆自動化"}
  @{R=8;  B="【急募】LINE WORKSで定期メッセージ配信ツール作成依頼";                                   C="システム開発"; D="5,000 円 ~ 10,000 円 / 固定";       E="期限情報なし"; F="https://www.lancers.jp/work/detail/5403166"; G=65;  H="◆ツール"}
  @{R=9;  B="【急募】教育系のWEBサイトの作成";                                              C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定";       E="期限情報なし"; F="https://www.lancers.jp/work/detail/5403527"; G=33;  H="◇サイト"}
  @{R=10; B="【ビジネスパートナー募集】英語が得意な方を探しています";                                    C="システム開発"; D="300,000 円 ~ 500,000 円 / 固定";     E="期限情報なし"; F="https://www.lancers.jp/work/detail/5403384"; G=25;  H=$null}
  @{R=11; B="初回 サブスクペイからCSVデータをダウンロードし、データベース同期するプログラムの作成";                     C="システム開発"; D="50,000 円 ~ 100,000 円 / 固定";      E="期限情報なし"; F="https://www.lancers.jp/work/detail/5403072"; G=18;  H=$null}
  @{R=12; B="限定公開 PR 限定公開の仕事";                                                C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定";       E="期限情報なし"; F="https://www.lancers.jp/work/detail/5399347"; G=13;  H=$null}
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $newTimestamp
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    if ($row.H -ne $null) {
        $ws.Cells.Item($r, 8).Value = $row.H
    }
}

# ---------------------------------------------------------------------
# 4) Re-create the column-F hyperlinks for every data row, in row order,
#    so the link target always matches the (possibly refreshed) URL text.
#    (This engine clears the whole sheet's hyperlink collection as soon
#    as any single one is deleted, so the simplest reliable path is to
#    drop them all once and rebuild top-to-bottom.)
# ---------------------------------------------------------------------
$ws.Range("F2").Hyperlinks.Delete()

foreach ($row in $rows) {
    $cell = $ws.Cells.Item($row.R, 6)
    $ws.Hyperlinks.Add($cell, $row.F)
    $cell.Style = "Hyperlink"
}
